$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Crew "Family name, given names" column (column 3) replacements.
# We look the name column up per-row and overwrite only the matching
# cell's Range.Text so that every other attribute of the surrounding
# run/paragraph (rsid, xml:space, etc.) is left completely untouched -
# a document-wide Find/Replace would regenerate the run and lose them.
$nameMap = @{
    "Garza, Daumantas Sawney"      = "Meyrick, BertrandLovise"
    "Rake, Batuhan Khodadad"       = "Genadiev, AyselKamil"
    "Evangelista, Veniamin Kyros"  = "Traylor, GovindaDiodotos"
    "Crespo, HonorataMarko"        = "Hakim, KalinLalita"
    "Furlan, IvkaWilla"            = "Post, Eun-JiIryna"
    "MacNevin, ApostolKanti"       = "Sanna, GuilhermeRamana"
    "Amador, DaudZinat"            = "Aitken, AdelaisEsdras"
    "Christian, VilhjálmurTerje"   = "Chiara, AniMaria"
    "Alessi, AuroraAmonet"         = "Wheatley, MelisaViraj"
    "Gwerder, MileJuliana"         = "Bridges, LauriLakshmana"
}

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $nameCell = $table.Cell($r, 3)
    $nameText = $nameCell.Range.Text
    foreach ($oldName in $nameMap.Keys) {
        if ($nameText -like "$oldName*") {
            $nameCell.Range.Text = $nameMap[$oldName]

            if ($oldName -eq "Gwerder, MileJuliana") {
                # The rank/rating cell ("UXO") belonging to this same row
                # must also change to "Client". "UXO" occurs several times
                # in the table (once per crew member with that rating), so
                # we only touch the cell in this specific row.
                $rankCell = $table.Cell($r, 5)
                if ($rankCell.Range.Text -like "UXO*") {
                    $rankCell.Range.Text = "Client"
                }
            }
            break
        }
    }
}
